$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row on the sheet (data currently lives in A:D)
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column before column B ("SFIA Level" and everything after
# shift one column to the right, so B:D becomes C:E) to hold the new
# "Skill Description" column.
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Populate the new column: for every data row, the skill's full/plain-text
# description mirrors the SkillCode already present in column A.
for ($r = 2; $r -le $lastRow; $r++) {
    $skillCode = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $skillCode
}
